$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC (31 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 540.1070999999999  # H28
$ws.Cells.Item(28, 9).Value = 540.1070999999999  # I28
$ws.Cells.Item(28, 11).Value = 540.1070999999999  # K28
$ws.Cells.Item(28, 13).Value = -55.10709999999995  # M28
$ws.Cells.Item(70, 8).Value = 2277.4  # H70
$ws.Cells.Item(70, 9).Value = 1999.5  # I70
$ws.Cells.Item(70, 11).Value = 5998.5  # K70
$ws.Cells.Item(70, 13).Value = -5728.5  # M70
$ws.Cells.Item(73, 8).Value = 2277.4  # H73
$ws.Cells.Item(73, 9).Value = 1999.5  # I73
$ws.Cells.Item(73, 11).Value = 5998.5  # K73
$ws.Cells.Item(73, 13).Value = -5062.5  # M73
$ws.Cells.Item(107, 8).Value = 47620956  # H107
$ws.Cells.Item(107, 9).Value = 55556650  # I107
$ws.Cells.Item(107, 11).Value = 55556650  # K107
$ws.Cells.Item(107, 13).Value = -55554730  # M107
$ws.Cells.Item(115, 9).Value = 111117144  # I115
$ws.Cells.Item(115, 11).Value = 333351432  # K115
$ws.Cells.Item(115, 13).Value = -333349865  # M115
$ws.Cells.Item(135, 8).Value = 14287440  # H135
$ws.Cells.Item(135, 9).Value = 17858486  # I135
$ws.Cells.Item(135, 11).Value = 160726374  # K135
$ws.Cells.Item(135, 13).Value = -160723839  # M135
$ws.Cells.Item(137, 8).Value = 3646.923  # H137
$ws.Cells.Item(137, 9).Value = 3178.3845  # I137
$ws.Cells.Item(137, 11).Value = 9535.1535  # K137
$ws.Cells.Item(137, 13).Value = -6985.1535  # M137
$ws.Cells.Item(138, 8).Value = 14497575  # H138
$ws.Cells.Item(138, 10).Value = 5364.5  # J138
$ws.Cells.Item(138, 12).Value = 16093.5  # L138
$ws.Cells.Item(138, 14).Value = -26373.5  # N138

# --- Worksheet: ARM (32 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2004.625  # H2
$ws.Cells.Item(2, 9).Value = 2004.9333  # I2
$ws.Cells.Item(2, 11).Value = 2004.9333  # K2
$ws.Cells.Item(2, 13).Value = -1891.9333  # M2
$ws.Cells.Item(32, 8).Value = 2348.4712  # H32
$ws.Cells.Item(32, 9).Value = 1758.0667  # I32
$ws.Cells.Item(32, 11).Value = 1758.0667  # K32
$ws.Cells.Item(32, 13).Value = -1471.0667  # M32
$ws.Cells.Item(61, 8).Value = 3812.35  # H61
$ws.Cells.Item(61, 9).Value = 3702.611  # I61
$ws.Cells.Item(61, 11).Value = 3702.611  # K61
$ws.Cells.Item(61, 13).Value = -3490.611  # M61
$ws.Cells.Item(97, 8).Value = 3862.8572  # H97
$ws.Cells.Item(97, 10).Value = 1753.6666  # J97
$ws.Cells.Item(97, 12).Value = 1753.6666  # L97
$ws.Cells.Item(97, 14).Value = -2745.6666  # N97
$ws.Cells.Item(102, 8).Value = 3203.532  # H102
$ws.Cells.Item(102, 9).Value = 2829.743  # I102
$ws.Cells.Item(102, 11).Value = 2829.743  # K102
$ws.Cells.Item(102, 13).Value = -1207.743  # M102
$ws.Cells.Item(116, 8).Value = 2004.625  # H116
$ws.Cells.Item(116, 9).Value = 2004.9333  # I116
$ws.Cells.Item(116, 11).Value = 2004.9333  # K116
$ws.Cells.Item(116, 13).Value = 289.0667000000001  # M116
$ws.Cells.Item(132, 8).Value = 6369.25  # H132
$ws.Cells.Item(132, 10).Value = 8406.362999999999  # J132
$ws.Cells.Item(132, 12).Value = 25219.089  # L132
$ws.Cells.Item(132, 14).Value = -30279.089  # N132
$ws.Cells.Item(136, 8).Value = 3812.35  # H136
$ws.Cells.Item(136, 9).Value = 3702.611  # I136
$ws.Cells.Item(136, 11).Value = 11107.833  # K136
$ws.Cells.Item(136, 13).Value = -8557.832999999999  # M136

# --- Worksheet: BSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2004.625  # H3
$ws.Cells.Item(3, 9).Value = 2004.9333  # I3
$ws.Cells.Item(3, 11).Value = 2004.9333  # K3
$ws.Cells.Item(3, 13).Value = -1890.9333  # M3
$ws.Cells.Item(22, 8).Value = 429.2  # H22
$ws.Cells.Item(22, 9).Value = 429.2  # I22
$ws.Cells.Item(22, 11).Value = 429.2  # K22
$ws.Cells.Item(22, 13).Value = -256.2  # M22
$ws.Cells.Item(86, 8).Value = 3085.5264  # H86
$ws.Cells.Item(86, 9).Value = 2827.7  # I86
$ws.Cells.Item(86, 10).Value = 3372  # J86
$ws.Cells.Item(86, 11).Value = 2827.7  # K86
$ws.Cells.Item(86, 12).Value = 3372  # L86
$ws.Cells.Item(86, 13).Value = -1704.7  # M86
$ws.Cells.Item(86, 14).Value = -5618  # N86
$ws.Cells.Item(89, 8).Value = 3085.5264  # H89
$ws.Cells.Item(89, 9).Value = 2827.7  # I89
$ws.Cells.Item(89, 10).Value = 3372  # J89
$ws.Cells.Item(89, 11).Value = 14138.5  # K89
$ws.Cells.Item(89, 12).Value = 16860  # L89
$ws.Cells.Item(89, 13).Value = -8522.5  # M89
$ws.Cells.Item(89, 14).Value = -28092  # N89

# --- Worksheet: CRP (47 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5075.227  # H31
$ws.Cells.Item(31, 10).Value = 5925.778  # J31
$ws.Cells.Item(31, 12).Value = 5925.778  # L31
$ws.Cells.Item(31, 14).Value = -6515.778  # N31
$ws.Cells.Item(34, 8).Value = 5075.227  # H34
$ws.Cells.Item(34, 10).Value = 5925.778  # J34
$ws.Cells.Item(34, 12).Value = 5925.778  # L34
$ws.Cells.Item(34, 14).Value = -6329.778  # N34
$ws.Cells.Item(58, 8).Value = 8057.643  # H58
$ws.Cells.Item(58, 9).Value = 9082.546  # I58
$ws.Cells.Item(58, 10).Value = 4299.6665  # J58
$ws.Cells.Item(58, 11).Value = 9082.546  # K58
$ws.Cells.Item(58, 12).Value = 4299.6665  # L58
$ws.Cells.Item(58, 13).Value = -8879.546  # M58
$ws.Cells.Item(58, 14).Value = -4705.6665  # N58
$ws.Cells.Item(59, 8).Value = 79115  # H59
$ws.Cells.Item(59, 10).Value = 79115  # J59
$ws.Cells.Item(59, 12).Value = 79115  # L59
$ws.Cells.Item(59, 14).Value = -81405  # N59
$ws.Cells.Item(94, 8).Value = 5867.25  # H94
$ws.Cells.Item(94, 9).Value = 4492.25  # I94
$ws.Cells.Item(94, 10).Value = 7242.25  # J94
$ws.Cells.Item(94, 11).Value = 4492.25  # K94
$ws.Cells.Item(94, 12).Value = 7242.25  # L94
$ws.Cells.Item(94, 13).Value = -4041.25  # M94
$ws.Cells.Item(94, 14).Value = -8144.25  # N94
$ws.Cells.Item(132, 8).Value = 8007.857  # H132
$ws.Cells.Item(132, 9).Value = 8393.210999999999  # I132
$ws.Cells.Item(132, 10).Value = 4347  # J132
$ws.Cells.Item(132, 11).Value = 25179.633  # K132
$ws.Cells.Item(132, 12).Value = 13041  # L132
$ws.Cells.Item(132, 13).Value = -22649.633  # M132
$ws.Cells.Item(132, 14).Value = -18101  # N132
$ws.Cells.Item(134, 8).Value = 5578.9165  # H134
$ws.Cells.Item(134, 9).Value = 4553.778  # I134
$ws.Cells.Item(134, 10).Value = 8654.333000000001  # J134
$ws.Cells.Item(134, 11).Value = 13661.334  # K134
$ws.Cells.Item(134, 12).Value = 25962.999  # L134
$ws.Cells.Item(134, 13).Value = -11126.334  # M134
$ws.Cells.Item(134, 14).Value = -31032.999  # N134
$ws.Cells.Item(136, 8).Value = 8057.643  # H136
$ws.Cells.Item(136, 9).Value = 9082.546  # I136
$ws.Cells.Item(136, 10).Value = 4299.6665  # J136
$ws.Cells.Item(136, 11).Value = 27247.638  # K136
$ws.Cells.Item(136, 12).Value = 12898.9995  # L136
$ws.Cells.Item(136, 13).Value = -24697.638  # M136
$ws.Cells.Item(136, 14).Value = -17998.9995  # N136

# --- Worksheet: CUL (15 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 114.04762  # H2
$ws.Cells.Item(2, 10).Value = 194.6  # J2
$ws.Cells.Item(2, 12).Value = 1167.6  # L2
$ws.Cells.Item(2, 14).Value = -1393.6  # N2
$ws.Cells.Item(122, 8).Value = 736  # H122
$ws.Cells.Item(122, 10).Value = 894.3  # J122
$ws.Cells.Item(122, 12).Value = 8048.7  # L122
$ws.Cells.Item(122, 14).Value = -12948.7  # N122
$ws.Cells.Item(136, 8).Value = 12243.833  # H136
$ws.Cells.Item(136, 9).Value = 12493.6  # I136
$ws.Cells.Item(136, 10).Value = 10995  # J136
$ws.Cells.Item(136, 11).Value = 37480.8  # K136
$ws.Cells.Item(136, 12).Value = 32985  # L136
$ws.Cells.Item(136, 13).Value = -32380.8  # M136
$ws.Cells.Item(136, 14).Value = -43185  # N136

# --- Worksheet: GSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5003485.5  # H80
$ws.Cells.Item(80, 9).Value = 14377422  # I80
$ws.Cells.Item(80, 10).Value = 4052.9333  # J80
$ws.Cells.Item(80, 11).Value = 14377422  # K80
$ws.Cells.Item(80, 12).Value = 4052.9333  # L80
$ws.Cells.Item(80, 13).Value = -14376424  # M80
$ws.Cells.Item(80, 14).Value = -6048.933300000001  # N80
$ws.Cells.Item(83, 8).Value = 5003485.5  # H83
$ws.Cells.Item(83, 9).Value = 14377422  # I83
$ws.Cells.Item(83, 10).Value = 4052.9333  # J83
$ws.Cells.Item(83, 11).Value = 71887110  # K83
$ws.Cells.Item(83, 12).Value = 20264.6665  # L83
$ws.Cells.Item(83, 13).Value = -71882118  # M83
$ws.Cells.Item(83, 14).Value = -30248.6665  # N83

# --- Worksheet: LTW (37 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3759.5715  # H7
$ws.Cells.Item(7, 9).Value = 4263.5  # I7
$ws.Cells.Item(7, 11).Value = 4263.5  # K7
$ws.Cells.Item(7, 13).Value = -4151.5  # M7
$ws.Cells.Item(100, 8).Value = 2949.9  # H100
$ws.Cells.Item(100, 9).Value = 2749  # I100
$ws.Cells.Item(100, 10).Value = 2972.2222  # J100
$ws.Cells.Item(100, 11).Value = 2749  # K100
$ws.Cells.Item(100, 12).Value = 2972.2222  # L100
$ws.Cells.Item(100, 13).Value = -2208  # M100
$ws.Cells.Item(100, 14).Value = -4054.2222  # N100
$ws.Cells.Item(122, 8).Value = 8086.5625  # H122
$ws.Cells.Item(122, 9).Value = 7556  # I122
$ws.Cells.Item(122, 11).Value = 22668  # K122
$ws.Cells.Item(122, 13).Value = -20218  # M122
$ws.Cells.Item(126, 8).Value = 3759.5715  # H126
$ws.Cells.Item(126, 9).Value = 4263.5  # I126
$ws.Cells.Item(126, 11).Value = 12790.5  # K126
$ws.Cells.Item(126, 13).Value = -10320.5  # M126
$ws.Cells.Item(132, 8).Value = 25423.908  # H132
$ws.Cells.Item(132, 9).Value = 28309.334  # I132
$ws.Cells.Item(132, 10).Value = 2917.6  # J132
$ws.Cells.Item(132, 11).Value = 84928.00199999999  # K132
$ws.Cells.Item(132, 12).Value = 8752.799999999999  # L132
$ws.Cells.Item(132, 13).Value = -82398.00199999999  # M132
$ws.Cells.Item(132, 14).Value = -13812.8  # N132
$ws.Cells.Item(133, 8).Value = 38000  # H133
$ws.Cells.Item(133, 10).Value = 38000  # J133
$ws.Cells.Item(133, 12).Value = 38000  # L133
$ws.Cells.Item(133, 14).Value = -43060  # N133
$ws.Cells.Item(136, 8).Value = 7778.68  # H136
$ws.Cells.Item(136, 9).Value = 3926.6365  # I136
$ws.Cells.Item(136, 10).Value = 10805.286  # J136
$ws.Cells.Item(136, 11).Value = 11779.9095  # K136
$ws.Cells.Item(136, 12).Value = 32415.858  # L136
$ws.Cells.Item(136, 13).Value = -9229.9095  # M136
$ws.Cells.Item(136, 14).Value = -37515.858  # N136

# --- Worksheet: WVR (26 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1233.5  # H113
$ws.Cells.Item(113, 9).Value = 1227.3438  # I113
$ws.Cells.Item(113, 11).Value = 3682.0314  # K113
$ws.Cells.Item(113, 13).Value = -1512.0314  # M113
$ws.Cells.Item(122, 8).Value = 13123.9375  # H122
$ws.Cells.Item(122, 9).Value = 6172  # I122
$ws.Cells.Item(122, 10).Value = 24710.5  # J122
$ws.Cells.Item(122, 11).Value = 18516  # K122
$ws.Cells.Item(122, 12).Value = 74131.5  # L122
$ws.Cells.Item(122, 13).Value = -16066  # M122
$ws.Cells.Item(122, 14).Value = -79031.5  # N122
$ws.Cells.Item(126, 8).Value = 4174.6665  # H126
$ws.Cells.Item(126, 9).Value = 3764.6  # I126
$ws.Cells.Item(126, 11).Value = 11293.8  # K126
$ws.Cells.Item(126, 13).Value = -8823.799999999999  # M126
$ws.Cells.Item(132, 8).Value = 3469  # H132
$ws.Cells.Item(132, 9).Value = 3392.9119  # I132
$ws.Cells.Item(132, 11).Value = 10178.7357  # K132
$ws.Cells.Item(132, 13).Value = -7648.735700000001  # M132
$ws.Cells.Item(136, 8).Value = 21112.309  # H136
$ws.Cells.Item(136, 9).Value = 23053.531  # I136
$ws.Cells.Item(136, 10).Value = 14900.4  # J136
$ws.Cells.Item(136, 11).Value = 69160.59299999999  # K136
$ws.Cells.Item(136, 12).Value = 44701.2  # L136
$ws.Cells.Item(136, 13).Value = -66610.59299999999  # M136
$ws.Cells.Item(136, 14).Value = -49801.2  # N136
